$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 28
$ws.Range("I6").Value = 13.6
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 40.8
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = 71.2
$ws.Range("N6").Value = -524
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H33").Value = 247.2
$ws.Range("I33").Value = 247.2
$ws.Range("K33").Value = 247.2
$ws.Range("M33").Value = -18.19999999999999
$ws.Range("H40").Value = 4594.2573
$ws.Range("J40").Value = 5327.5864
$ws.Range("L40").Value = 5327.5864
$ws.Range("N40").Value = -5677.5864
$ws.Range("H43").Value = 7285.7144
$ws.Range("J43").Value = 9750
$ws.Range("L43").Value = 9750
$ws.Range("N43").Value = -9888
$ws.Range("H55").Value = 234.28572
$ws.Range("I55").Value = 128
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 128
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = 86
$ws.Range("N55").Value = -928
$ws.Range("H137").Value = 1843.6
$ws.Range("I137").Value = 1843.6
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5530.799999999999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -2980.799999999999
$ws.Range("N137").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H17").Value = 2475
$ws.Range("I17").Value = 2450
$ws.Range("J17").Value = 2500
$ws.Range("K17").Value = 2450
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = -2277
$ws.Range("N17").Value = -2846
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H21").Value = 500
$ws.Range("I21").Value = 500
$ws.Range("K21").Value = 500
$ws.Range("M21").Value = -126
$ws.Range("H22").Value = 2508
$ws.Range("I22").Value = 2016
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 2016
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -1717
$ws.Range("N22").Value = -3598
$ws.Range("H23").Value = 16000
$ws.Range("J23").Value = 16000
$ws.Range("L23").Value = 16000
$ws.Range("N23").Value = -16518
$ws.Range("H25").Value = 715.8
$ws.Range("I25").Value = 715.8
$ws.Range("K25").Value = 715.8
$ws.Range("M25").Value = -313.8
$ws.Range("H30").Value = 498.75
$ws.Range("I30").Value = 498.75
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 498.75
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -348.75
$ws.Range("N30").ClearContents()
$ws.Range("H32").Value = 4477.4165
$ws.Range("I32").Value = 4477.4165
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4477.4165
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4190.4165
$ws.Range("N32").ClearContents()
$ws.Range("H132").Value = 6752
$ws.Range("I132").Value = 6353.8
$ws.Range("J132").Value = 7747.5
$ws.Range("K132").Value = 19061.4
$ws.Range("L132").Value = 23242.5
$ws.Range("M132").Value = -16531.4
$ws.Range("N132").Value = -28302.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H22").Value = 4075
$ws.Range("I22").Value = 4100
$ws.Range("K22").Value = 4100
$ws.Range("M22").Value = -3927
$ws.Range("H86").Value = 1071
$ws.Range("I86").Value = 1142.5
$ws.Range("J86").Value = 785
$ws.Range("K86").Value = 1142.5
$ws.Range("L86").Value = 785
$ws.Range("M86").Value = -19.5
$ws.Range("N86").Value = -3031
$ws.Range("H89").Value = 1071
$ws.Range("I89").Value = 1142.5
$ws.Range("J89").Value = 785
$ws.Range("K89").Value = 5712.5
$ws.Range("L89").Value = 3925
$ws.Range("M89").Value = -96.5
$ws.Range("N89").Value = -15157
$ws.Range("H106").Value = 26625
$ws.Range("J106").Value = 26625
$ws.Range("L106").Value = 26625
$ws.Range("N106").Value = -29149
$ws.Range("H134").Value = 3989.8667
$ws.Range("I134").Value = 1635.3
$ws.Range("K134").Value = 4905.9
$ws.Range("M134").Value = -2370.9

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 4211.75
$ws.Range("I36").Value = 4211.75
$ws.Range("K36").Value = 4211.75
$ws.Range("M36").Value = -3823.75
$ws.Range("H40").Value = 4211.75
$ws.Range("I40").Value = 4211.75
$ws.Range("K40").Value = 4211.75
$ws.Range("M40").Value = -4051.75
$ws.Range("H50").Value = 24078.934
$ws.Range("I50").Value = 30000
$ws.Range("J50").Value = 23168
$ws.Range("K50").Value = 30000
$ws.Range("L50").Value = 23168
$ws.Range("M50").Value = -29375
$ws.Range("N50").Value = -24418
$ws.Range("H51").Value = 17954.545
$ws.Range("J51").Value = 17954.545
$ws.Range("L51").Value = 17954.545
$ws.Range("N51").Value = -19426.545
$ws.Range("H60").Value = 23636.363
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 25000
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 25000
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -26022
$ws.Range("H61").Value = 17954.545
$ws.Range("J61").Value = 17954.545
$ws.Range("L61").Value = 17954.545
$ws.Range("N61").Value = -18650.545
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value = 778.2
$ws.Range("J107").Value = 796
$ws.Range("L107").Value = 796
$ws.Range("N107").Value = -4636
$ws.Range("H122").Value = 1778.8
$ws.Range("I122").Value = 1723.5
$ws.Range("K122").Value = 5170.5
$ws.Range("M122").Value = -2720.5
$ws.Range("H134").Value = 9874.25
$ws.Range("I134").Value = 9832.666999999999
$ws.Range("K134").Value = 29498.001
$ws.Range("M134").Value = -26963.001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50.31579
$ws.Range("I2").Value = 18.083334
$ws.Range("K2").Value = 108.500004
$ws.Range("M2").Value = 4.499995999999996
$ws.Range("H14").Value = 795.1667
$ws.Range("I14").Value = 795.1667
$ws.Range("K14").Value = 2385.5001
$ws.Range("M14").Value = -2212.5001
$ws.Range("H39").Value = 4305.5
$ws.Range("I39").Value = 4305.5
$ws.Range("K39").Value = 12916.5
$ws.Range("M39").Value = -12622.5
$ws.Range("H92").Value = 2866.6667
$ws.Range("I92").Value = 2550
$ws.Range("J92").Value = 3500
$ws.Range("K92").Value = 7650
$ws.Range("L92").Value = 10500
$ws.Range("M92").Value = -6402
$ws.Range("N92").Value = -12996
$ws.Range("H117").Value = 917.5714
$ws.Range("I117").Value = 607.5
$ws.Range("J117").Value = 1331
$ws.Range("K117").Value = 1822.5
$ws.Range("L117").Value = 3993
$ws.Range("M117").Value = 1619.5
$ws.Range("N117").Value = -10877
$ws.Range("H121").Value = 1156
$ws.Range("I121").Value = 731.8333
$ws.Range("J121").Value = 1474.125
$ws.Range("K121").Value = 2195.4999
$ws.Range("L121").Value = 4422.375
$ws.Range("M121").Value = -885.4998999999998
$ws.Range("N121").Value = -7042.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 21.666666
$ws.Range("I2").Value = 23
$ws.Range("K2").Value = 23
$ws.Range("M2").Value = 90
$ws.Range("H57").Value = 28750
$ws.Range("J57").Value = 28750
$ws.Range("L57").Value = 28750
$ws.Range("N57").Value = -30390

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2110.8
$ws.Range("I16").Value = 2110.8
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2110.8
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1940.8
$ws.Range("N16").ClearContents()
$ws.Range("H32").Value = 1766.3334
$ws.Range("I32").Value = 1766.3334
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1766.3334
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1449.3334
$ws.Range("N32").ClearContents()
$ws.Range("H93").Value = 7959.875
$ws.Range("I93").Value = 7959.875
$ws.Range("K93").Value = 7959.875
$ws.Range("M93").Value = -6711.875
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 1831.6666
$ws.Range("I136").Value = 1831.6666
$ws.Range("K136").Value = 5494.9998
$ws.Range("M136").Value = -2944.9998

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 19767
$ws.Range("J80").Value = 19767
$ws.Range("L80").Value = 19767
$ws.Range("N80").Value = -21763
$ws.Range("H83").Value = 19767
$ws.Range("J83").Value = 19767
$ws.Range("L83").Value = 59301
$ws.Range("N83").Value = -69285
$ws.Range("H126").Value = 1905.3334
$ws.Range("I126").Value = 1905.3334
$ws.Range("K126").Value = 5716.0002
$ws.Range("M126").Value = -3246.0002
$ws.Range("H132").Value = 8495.5
$ws.Range("I132").Value = 8055.222
$ws.Range("K132").Value = 24165.666
$ws.Range("M132").Value = -21635.666
$ws.Range("H136").Value = 1742.5714
$ws.Range("I136").Value = 1742.5714
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5227.7142
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2677.7142
$ws.Range("N136").ClearContents()
